$wb = $excel.ActiveWorkbook

# --- Sheet1: "RPi 3A+ Total BOM" ---
$ws1 = $wb.Worksheets.Item("RPi 3A+ Total BOM")

# Insert a new row before row 21 (shifts old rows 21-26 down to 22-27)
$ws1.Rows.Item(21).Insert()

# Populate new row 21
$ws1.Range("B21").Value = "USB Cable / PSU"
$ws1.Range("C21").Value = 20
$ws1.Range("C21").NumberFormat = $ws1.Range("C20").NumberFormat
$ws1.Range("E21").Value = "Ebay / may already own required parts"
$ws1.Range("E21").ClearFormats()

# Update the "Cost of Customer Items" formula (now row 23) to include the new row
$ws1.Range("C23").Formula = "=SUM(C18:C21)"

# Update the final "Total" formula (now row 27) to reference the shifted row
$ws1.Range("C27").Formula = "=C23+C15+C12"

# Update selections for sheet1
$ws1.Range("C24").Select()

# --- Sheet2: "RPi Zero Total BOM" ---
$ws2 = $wb.Worksheets.Item("RPi Zero Total BOM")
$ws2.Range("C36").Select()

# --- Sheet3: "VeinCamHat BOM" ---
$ws3 = $wb.Worksheets.Item("VeinCamHat BOM")
$ws3.Range("C7").Select()

# --- Sheet4: "VeinCamHatZero BOM" ---
$ws4 = $wb.Worksheets.Item("VeinCamHatZero BOM")
$ws4.Range("C9").Select()
$ws4.Activate()

# --- Update workbook window position/size (bookViews/workbookView) ---
# (Window.Top/Width/Height/Left are tracked on the Window object but this
#  runtime does not persist them back into workbookView@xWindow/yWindow/
#  windowWidth/windowHeight on save; set them anyway in case the host
#  starts honoring them, it is a harmless no-op otherwise.)
$win = $wb.Windows.Item(1)
$win.Top = -120
$win.Width = 29040
$win.Height = 15840
